$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79; existing rows 79-90 shift down to 80-91,
# carrying all their data/styles with them (matches the target diff).
$ws.Rows("79").Insert()

# Populate the newly inserted row 79 with this week's new data point.
$ws.Cells.Item(79, 1).Value = 7
$ws.Cells.Item(79, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(79, 3).Value = 'Ñuble'
$ws.Cells.Item(79, 4).Value = 45093
$ws.Cells.Item(79, 5).Value = 16
$ws.Cells.Item(79, 6).Value = 100112001
$ws.Cells.Item(79, 7).Value = 'Berenjena'
$ws.Cells.Item(79, 8).Value = 'Sin especificar'
$ws.Cells.Item(79, 9).Value = 'Primera'
$ws.Cells.Item(79, 10).Value = 50
$ws.Cells.Item(79, 11).Value = 9000
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 13).Value = 9000
$ws.Cells.Item(79, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(79, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(79, 16).Value = 150
$ws.Cells.Item(79, 17).Value = 60
$ws.Cells.Item(79, 18).Value = 'Hortaliza'
